$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..27 (col A = numeric id, col B = ONREG ticket string)
$data = @(
    @(2, 100147025, "ONREG-23902"),
    @(3, 101487374, "ONREG-17938"),
    @(4, 102440213, "ONREG-22294"),
    @(5, 110062018, "ONREG-18131"),
    @(6, 127178114, "ONREG-18240"),
    @(7, 139159854, "ONREG-25876"),
    @(8, 140389560, "ONREG-25888"),
    @(9, 102440213, "ONREG-22294"),
    @(10, 137199660, "ONREG-19423"),
    @(11, 100581355, "ONREG-19642"),
    @(12, 135142501, "ONREG-19769"),
    @(13, 102522400, "ONREG-19795"),
    @(14, 109514216, "ONREG-19814"),
    @(15, 119002359, "ONREG-19836"),
    @(16, 108482050, "ONREG-18965"),
    @(17, 163203979, "ONREG-23937"),
    @(18, 146457800, "ONREG-20446"),
    @(19, 119592056, "ONREG-20191"),
    @(20, 164660957, "ONREG-12902"),
    @(21, 113526468, "ONREG-12961"),
    @(22, 150744017, "ONREG-13268"),
    @(23, 131704213, "ONREG-13351"),
    @(24, 131478014, "ONREG-13415"),
    @(25, 139637466, "ONREG-12927"),
    @(26, 142320256, "ONREG-13024"),
    @(27, 124700667, "ONREG-13036")
)

# Cells beyond the sheet's original used range default to Text on write
# unless a NumberFormat is established first -- prime the whole target
# range so numeric values land as real numbers.
$ws.Range("A2:B27").NumberFormat = "General"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

# Remove the custom row height (15.75) that used to mark the data rows —
# rows are back to the sheet default height.
$ws.Range("A2:B27").RowHeight = 15

# Strip the old per-row font formatting (Segoe UI / Times New Roman) so the
# table goes back to the workbook's default font, keeping A2 bold.
$ws.Range("A2:B27").Font.Name = "Calibri"
$ws.Range("A2:B27").Font.Size = 11
$ws.Range("A2:B27").Font.Bold = $false
$ws.Range("A2").Font.Bold = $true

$ws.Range("A18").Select()
